$d = $word.ActiveDocument

# 1. Insert two new "Team Number:" / "Team Member Names:" paragraphs right
#    after the "Keep answers as short as possible..." paragraph (and before
#    the "Provide acceleration plots and analysis" heading).
$anchorRange = $d.Content
$anchorRange.Find.Execute("Keep answers as short as possible while still meeting specifications. Submit as a PDF.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $anchorRange.Paragraphs(1)

$insertRange = $anchorPara.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()

$newPara1 = $anchorPara.Next()
$newPara1.Range.Text = "Team Number:"
$newPara1.Style = "BodyText"

$afterRange = $newPara1.Range
$afterRange.Collapse(0)
$afterRange.InsertParagraphAfter()

$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "Team Member Names:"
$newPara2.Style = "BodyText"

# 2. Update checklist item wording.
$d.Content.Find.Execute("Calcuations that relate acceleration due to gravity to one Teensy unit are correct.", $true, $false, $false, $false, $false, $true, 1, $false, "Calcuations that relate acceleration due to gravity to one accelerometer unit are correct.", 2)

$d.Content.Find.Execute("Take the mean of data to calculate the resting zero acceleration values. Report them clearly.", $true, $false, $false, $false, $false, $true, 1, $false, "Take the mean of data to calculate the resting zero acceleration values. Report them clearly with appropriate uncertainty bounds.", 2)

$d.Content.Find.Execute("Includes description of how to calculate descriptive statistics and correct values.", $true, $false, $false, $false, $false, $true, 1, $false, "Includes description of how to calculate statistics and correct values.", 2)
